$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.387.94"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.852.24"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'240.59"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'0.6341"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.07585"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'24.50"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.850.12"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'0.6803"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "'83.34"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'6.133"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "29.363.95"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'230.14"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'7.467"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'158.66"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'17.66"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "'1.410"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "'1.476"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "'0.05684"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'4.127"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'4.051"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "'1.829"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "'1.156"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "'0.6996"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'2.577"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "1.245.17"
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("D39").Value = "'2.727"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "'6.421"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "'0.9051"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'0.9990"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "2.007.89"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "'102.37"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "'65.86"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'7.139"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'0.1167"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "'9.033"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'0.00000000115"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "'0.3959"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("E51").Value = "  -0.15%  "
